$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1:AF1 with the same style as the existing header row
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season-record values (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 63; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}
